# Apply the edit described by the diff:
#  - Row 7's time slot becomes "19:0-19:5" (previously "22:30-22:35")
#    and its "Посыл / Заповедь" cell (B7) switches to the same message
#    used in rows 3/5 (instead of the one used in rows 2/4/6).
#  - Rows 8-11 (the trailing 22:35-22:40 .. 22:50-22:55 entries) are removed.
#  - The active selection moves from B15 to B13 to reflect the now-shorter sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last four data rows (old rows 8-11); this also shifts the
# sheet's used range / dimension up from A1:C11 to A1:C7 automatically.
$ws.Rows("8:11").Delete()

# Row 7 now becomes the final data row: update its message/time pair.
$ws.Range("B7").Value2 = $ws.Range("B3").Value2
$ws.Range("C7").Value2 = "19:0-19:5"

# Reflect the new, shorter extent in the sheet's selection.
$ws.Range("B13").Select()
